$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet
$ws.Name = "Fed Poverty Calculations 2018"

# Update cell T1 from the number 2019 to the text "no"
$ws.Range("T1").Value = "no"

# Hide columns D, F and I (previously visible / bestFit)
$ws.Columns("D").Hidden = $true
$ws.Columns("F").Hidden = $true
$ws.Columns("I").Hidden = $true

# Move the active selection to G10
$ws.Range("G10").Select() | Out-Null
